$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the existing accelerometer samples (rows 2-22) before shifting them.
$data = @()
for ($r = 2; $r -le 22; $r++) {
    $x = $ws.Cells.Item($r,1).Value2
    $y = $ws.Cells.Item($r,2).Value2
    $z = $ws.Cells.Item($r,3).Value2
    $data += ,@($x,$y,$z)
}

# Two new readings are inserted at the start of the series.
$ws.Cells.Item(2,1).Value = -1.70257568359375
$ws.Cells.Item(2,2).Value = 5.600172996520996
$ws.Cells.Item(2,3).Value = 2.819403648376465

$ws.Cells.Item(3,1).Value = -2.645308017730713
$ws.Cells.Item(3,2).Value = 5.044093132019043
$ws.Cells.Item(3,3).Value = 3.529041290283203

# The previously captured rows now shift down two rows. Only the first
# 18 of the 21 captured rows are kept (old rows 2-19), which together
# with the two new rows above gives 20 data rows total; the last three
# captured rows (old rows 20-22) are dropped, net one row shorter overall.
$keep = $data.Length - 3
for ($i = 0; $i -lt $keep; $i++) {
    $r = $i + 4
    $ws.Cells.Item($r,1).Value = $data[$i][0]
    $ws.Cells.Item($r,2).Value = $data[$i][1]
    $ws.Cells.Item($r,3).Value = $data[$i][2]
}

# Remove the now-unused trailing row left over from the old layout.
$ws.Rows.Item(22).Delete()

$ws.Range("A1").Select()
